$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Obrigatorio" flag (column E) from "N" to "S" for rows 2-8
$ws.Range("E2:E8").Value = "S"
